# Adds and drafts out the project template page:
# introduces a new "case study" column (E) with a boolean value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + value for column E
$ws.Range("E1").Value = "case study"
$ws.Range("E2").Value = $false

# Column widths (A-D now get explicit custom widths matching the target layout)
$ws.Columns.Item(1).ColumnWidth = 11.5
$ws.Columns.Item(2).ColumnWidth = 12.166666666666666
$ws.Columns.Item(3).ColumnWidth = 12.333333333333334
$ws.Columns.Item(4).ColumnWidth = 11

# Move the selection to the newly added cell
$ws.Range("E2").Select()
